$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Event_data")
$ws1.Range("B2").Value = 28444292
$ws1.Range("B3").Value = 22247864
$ws1.Range("B4").Value = 2508905
$ws1.Range("B5").Value = 937804
$ws1.Range("B6").Value = 211989
$ws1.Range("B7").Value = 171674
$ws1.Range("B8").Value = 2114127
$ws1.Range("B9").Value = 5324879
$ws1.Range("B10").Value = 3362109
$ws1.Range("B11").Value = 320394
$ws1.Range("B12").Value = 140506
$ws1.Range("B13").Value = 60907472
$ws1.Range("B14").Value = 97798729
$ws1.Range("B15").Value = 129939633
$ws1.Range("B16").Value = 156184630
$ws1.Range("B17").Value = 118683194
$ws1.Range("B18").Value = 72741773
$ws1.Range("B19").Value = 140404
$ws1.Range("B20").Value = 44314
$ws1.Range("B21").Value = 18803
$ws1.Range("B22").Value = 114128
$ws1.Range("B23").Value = 25864
$ws1.Range("B24").Value = 11299
$ws1.Range("B25").Value = 4279
$ws1.Range("B26").Value = 176586
$ws1.Range("B27").Value = 5320
$ws1.Range("B28").Value = 10390
$ws1.Range("B29").Value = 5900

$ws2 = $wb.Worksheets.Item("Energy_data")
$ws2.Range("B2").Value = 358398079.2
$ws2.Range("B3").Value = 323038985.28
$ws2.Range("B4").Value = 38386246.5
$ws2.Range("B5").Value = 17143057.12
$ws2.Range("B6").Value = 4260978.9
$ws2.Range("B7").Value = 3505583.08
$ws2.Range("B8").Value = 41584878.09
$ws2.Range("B9").Value = 24920433.72
$ws2.Range("B10").Value = 16642439.55
$ws2.Range("B11").Value = 3030927.24
$ws2.Range("B12").Value = 1296114.280096082
$ws2.Range("B13").Value = 22048504.864
$ws2.Range("B14").Value = 18581758.51
$ws2.Range("B15").Value = 48597422.742
$ws2.Range("B16").Value = 25301910.06
$ws2.Range("B17").Value = 925728.9132
$ws2.Range("B18").Value = 945643.049
$ws2.Range("B19").Value = 1432120.8
$ws2.Range("B20").Value = 535313.12
$ws2.Range("B21").Value = 239738.25
$ws2.Range("B22").Value = 215701.92
$ws2.Range("B23").Value = 65953.2
$ws2.Range("B24").Value = 32315.14
$ws2.Range("B25").Value = 12965.37
$ws2.Range("B26").Value = 508567.68
$ws2.Range("B27").Value = 34473.60000000001
$ws2.Range("B28").Value = 77821.10000000001
$ws2.Range("B29").Value = 47200

$ws3 = $wb.Worksheets.Item("Species_data")
$ws3.Range("B2").Value = 320394
$ws3.Range("B3").Value = 211989
$ws3.Range("B4").Value = 3533783
$ws3.Range("B5").Value = 2508905
$ws3.Range("B6").Value = 6403189
$ws3.Range("B7").Value = 22247864
$ws3.Range("B8").Value = 28444292
$ws3.Range("B9").Value = 2114127
$ws3.Range("B10").Value = 28105126
$ws3.Range("B11").Value = 937804
$ws3.Range("B12").Value = 140506
$ws3.Range("B13").Value = 10631651
$ws3.Range("B14").Value = 171674

Write-Host "Updated all values"